# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    populated with the fund holdings for that quarter.
# 2. Insert a new leading data row into the "总计" sheet summarizing the
#    2022-Q1 quarter, shifting the existing rows down and renumbering the
#    index column.

function Set-TextValue($cell, [string]$text) {
    # Force the cell to keep its value as text (even when it looks numeric,
    # e.g. fund codes like "001364" or decimal-looking figures like "4.72"),
    # then drop back to the default ("Normal") style so no extra / stray
    # cell formatting (like the quote-prefix style) is left behind.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Seed the new sheet from an existing quarter sheet so it inherits the exact
# same header/row styling (bordered bold header + bold index column).
$templateSheet = $wb.Worksheets.Item("2021-Q3")
$templateSheet.Range("A1:H2").Copy($newSheet.Range("A1:H2"))

# Replicate the data-row formatting down for the remaining 4 fund rows.
$dataRowTemplate = $newSheet.Range("A2:H2")
for ($r = 3; $r -le 6; $r++) {
    $dataRowTemplate.Copy($newSheet.Range("A" + $r + ":H" + $r))
}

# The template's A1 is blank; drop the stray empty cell the copy introduces.
$newSheet.Cells.Item(1, 1).ClearContents()

# Header row (plain text, not numeric-looking, so no quote-prefix trick /
# style reset is needed here -- keep the bold bordered header style intact).
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

$fundRows = @(
    @{ Idx = 0; Code = "001364"; Name = "大成景润灵活配置混合";     Size = "4.72"; Pos = "26.73"; Pct = "2.27"; Value = "0.1071"; Rank = 3 },
    @{ Idx = 1; Code = "002383"; Name = "大成趋势回报灵活配置混合"; Size = "4.64"; Pos = "24.42"; Pct = "1.99"; Value = "0.0923"; Rank = 3 },
    @{ Idx = 2; Code = "003373"; Name = "大成景禄灵活配置混合A";    Size = "3.83"; Pos = "28.91"; Pct = "2.30"; Value = "0.0881"; Rank = 5 },
    @{ Idx = 3; Code = "003374"; Name = "大成景禄灵活配置混合C";    Size = "1.93"; Pos = "28.91"; Pct = "2.30"; Value = "0.0444"; Rank = 5 },
    @{ Idx = 4; Code = "003147"; Name = "大成动态量化配置策略混合"; Size = "1.50"; Pos = "27.36"; Pct = "2.01"; Value = "0.0302"; Rank = 7 }
)

$row = 2
foreach ($f in $fundRows) {
    $newSheet.Cells.Item($row, 1).Value = $f.Idx
    Set-TextValue $newSheet.Cells.Item($row, 2) $f.Code
    Set-TextValue $newSheet.Cells.Item($row, 3) $f.Name
    Set-TextValue $newSheet.Cells.Item($row, 4) $f.Size
    Set-TextValue $newSheet.Cells.Item($row, 5) $f.Pos
    Set-TextValue $newSheet.Cells.Item($row, 6) $f.Pct
    Set-TextValue $newSheet.Cells.Item($row, 7) $f.Value
    $newSheet.Cells.Item($row, 8).Value = $f.Rank
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Shift the existing rows (2..4) down to (3..5), copying values + formatting,
# bottom-up so nothing gets overwritten prematurely.
for ($r = 4; $r -ge 2; $r--) {
    $newR = $r + 1
    $totalSheet.Range("A" + $r + ":D" + $r).Copy($totalSheet.Range("A" + $newR + ":D" + $newR))
}

# Write the new 2022-Q1 summary into row 2.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 5
$totalSheet.Cells.Item(2, 4).Value = 0.36

# Renumber the index column for the shifted rows.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
